$wb = $excel.ActiveWorkbook

# Update the "zh-cn" sheet: row 16 handoff/handback datetimes
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D16").Value = "2016-03-08 06:56:40"
$wsZhCn.Range("G16").Value = "2016-03-08 06:57:22"

# Update the "de-de" sheet: row 16 handoff/handback datetimes
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D16").Value = "2016-03-08 06:56:50"
$wsDeDe.Range("G16").Value = "2016-03-08 06:57:37"
